$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instruction Set")
$ws.Activate()

# Swap the Operand columns (E/F) for the shift-instruction rows: the
# "<reg>" / "<nnnn>[!<reg>]" values were in the wrong columns.
$ws.Range("E11").Value = "<reg>"
$ws.Range("F11").Value = "<nnnn>[!<reg>]"

$ws.Range("E12").Value = "<reg>"
$ws.Range("F12").Value = "<nnnn>[!<reg>]"

$ws.Range("E13").Value = "<reg>"
$ws.Range("F13").Value = "<nnnn>[!<reg>]"

$ws.Range("E14").Value = "<reg>"
$ws.Range("F14").Value = "<nnnn>[!<reg>]"

# jpovr (row 23) has no second operand - mark it with a literal "-"
# (quote-prefixed so Excel treats it as text, not a formula/negative).
$ws.Range("E23").Value = "'-"

# jp (row 26): the second operand was misplaced in E; it belongs in F,
# and E should show the "-" placeholder like the other no-operand rows.
$ws.Range("E26").Value = "'-"
$ws.Range("F26").Value = "<nnnn>[!<reg>]"

# Make "Instruction Set" the selected/active tab with E30 selected,
# instead of "Standard Programming Practice".
$ws.Range("E30").Select()
